$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted above the existing row 215,
# pushing the whole block (old rows 215-291) down by one row.
$ws.Rows(215).Insert()

# New row 215 mirrors the (now shifted) row 216 data, except for an
# updated date and updated min/max/avg/kg prices (a new week's reading).
$ws.Range("A215").Value = 7
$ws.Range("B215").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C215").Value = "Ñuble"
$ws.Range("D215").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D215").Value = 44468
$ws.Range("E215").Value = 16
$ws.Range("F215").Value = 100112020
$ws.Range("G215").Value = "Tomate"
$ws.Range("H215").Value = "Larga vida"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 160
$ws.Range("K215").Value = 6000
$ws.Range("L215").Value = 6500
$ws.Range("M215").Value = 6250
$ws.Range("N215").Value = "$/caja 10 kilos"
$ws.Range("O215").Value = "Región de Arica y Parinacota"
$ws.Range("P215").Value = 625
$ws.Range("Q215").Value = 10
$ws.Range("R215").Value = "Hortaliza"
